$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 45294.5
$ws.Range("I68").Value = 40294
$ws.Range("K68").Value = 40294
$ws.Range("M68").Value = -39545
$ws.Range("H71").Value = 45294.5
$ws.Range("I71").Value = 40294
$ws.Range("K71").Value = 120882
$ws.Range("M71").Value = -117138
$ws.Range("H76").Value = 1699
$ws.Range("I76").Value = 1699
$ws.Range("K76").Value = 1699
$ws.Range("M76").Value = -1384
$ws.Range("H79").Value = 1699
$ws.Range("I79").Value = 1699
$ws.Range("K79").Value = 1699
$ws.Range("M79").Value = -607
$ws.Range("H99").Value = 6868.5713
$ws.Range("I99").Value = 200
$ws.Range("J99").Value = 7980
$ws.Range("K99").Value = 600
$ws.Range("L99").Value = 23940
$ws.Range("M99").Value = 898
$ws.Range("N99").Value = -26936
$ws.Range("H100").Value = 2356.7144
$ws.Range("I100").Value = 1900.2
$ws.Range("J100").Value = 3498
$ws.Range("K100").Value = 1900.2
$ws.Range("L100").Value = 3498
$ws.Range("M100").Value = -1359.2
$ws.Range("N100").Value = -4580
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H103").Value = 12671.375
$ws.Range("I103").Value = 20537.25
$ws.Range("K103").Value = 61611.75
$ws.Range("M103").Value = -61025.75
$ws.Range("H112").Value = 1756.0769
$ws.Range("J112").Value = 1369.8889
$ws.Range("L112").Value = 4109.6667
$ws.Range("N112").Value = -6325.6667
$ws.Range("H113").Value = 4001.25
$ws.Range("I113").Value = 4001.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4001.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -747.25
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 20893.072
$ws.Range("I132").Value = 21450.3
$ws.Range("K132").Value = 64350.89999999999
$ws.Range("M132").Value = -61820.89999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -330
$ws.Range("N10").ClearContents()
$ws.Range("H31").Value = 52499.75
$ws.Range("I31").Value = 52499.75
$ws.Range("K31").Value = 52499.75
$ws.Range("M31").Value = -52205.75
$ws.Range("H61").Value = 2012.3158
$ws.Range("I61").Value = 1337.5
$ws.Range("K61").Value = 1337.5
$ws.Range("M61").Value = -1125.5
$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26248
$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 25000
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81240
$ws.Range("H102").Value = 8933515
$ws.Range("I102").Value = 15626486
$ws.Range("J102").Value = 9553.166999999999
$ws.Range("K102").Value = 15626486
$ws.Range("L102").Value = 9553.166999999999
$ws.Range("M102").Value = -15624864
$ws.Range("N102").Value = -12797.167
$ws.Range("H136").Value = 2012.3158
$ws.Range("I136").Value = 1337.5
$ws.Range("K136").Value = 4012.5
$ws.Range("M136").Value = -1462.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3249.25
$ws.Range("J20").Value = 3766.3333
$ws.Range("L20").Value = 3766.3333
$ws.Range("N20").Value = -4260.3333
$ws.Range("H80").Value = 887.5454999999999
$ws.Range("I80").Value = 418.2
$ws.Range("K80").Value = 418.2
$ws.Range("M80").Value = 579.8
$ws.Range("H83").Value = 887.5454999999999
$ws.Range("I83").Value = 418.2
$ws.Range("K83").Value = 2091
$ws.Range("M83").Value = 2901
$ws.Range("H86").Value = 6892.375
$ws.Range("I86").Value = 2380
$ws.Range("J86").Value = 9599.799999999999
$ws.Range("K86").Value = 2380
$ws.Range("L86").Value = 9599.799999999999
$ws.Range("M86").Value = -1257
$ws.Range("N86").Value = -11845.8
$ws.Range("H89").Value = 6892.375
$ws.Range("I89").Value = 2380
$ws.Range("J89").Value = 9599.799999999999
$ws.Range("K89").Value = 11900
$ws.Range("L89").Value = 47999
$ws.Range("M89").Value = -6284
$ws.Range("N89").Value = -59231
$ws.Range("H94").Value = 1090.2
$ws.Range("I94").Value = 1090.2
$ws.Range("K94").Value = 1090.2
$ws.Range("M94").Value = -639.2
$ws.Range("I102").Value = 35000
$ws.Range("K102").Value = 35000
$ws.Range("M102").Value = -31755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 29161.75
$ws.Range("J43").Value = 29161.75
$ws.Range("L43").Value = 29161.75
$ws.Range("N43").Value = -29529.75
$ws.Range("H60").Value = 40144.668
$ws.Range("I60").Value = 9000
$ws.Range("J60").Value = 79075.5
$ws.Range("K60").Value = 9000
$ws.Range("L60").Value = 79075.5
$ws.Range("M60").Value = -8489
$ws.Range("N60").Value = -80097.5
$ws.Range("H68").Value = 78882.5
$ws.Range("J68").Value = 78882.5
$ws.Range("L68").Value = 78882.5
$ws.Range("N68").Value = -80380.5
$ws.Range("H71").Value = 78882.5
$ws.Range("J71").Value = 78882.5
$ws.Range("L71").Value = 236647.5
$ws.Range("N71").Value = -244135.5
$ws.Range("H101").Value = 29161.75
$ws.Range("J101").Value = 29161.75
$ws.Range("L101").Value = 29161.75
$ws.Range("N101").Value = -35651.75
$ws.Range("H106").Value = 22933.334
$ws.Range("J106").Value = 22933.334
$ws.Range("L106").Value = 22933.334
$ws.Range("N106").Value = -25457.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1957.8
$ws.Range("I129").Value = 1933
$ws.Range("K129").Value = 5799
$ws.Range("M129").Value = -799

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 17787.5
$ws.Range("J40").Value = 17787.5
$ws.Range("L40").Value = 17787.5
$ws.Range("N40").Value = -18089.5
$ws.Range("H43").Value = 19033.166
$ws.Range("J43").Value = 19839.8
$ws.Range("L43").Value = 19839.8
$ws.Range("N43").Value = -20141.8
$ws.Range("H55").Value = 4757.5
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H62").Value = 50000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 50000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 50000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -51372
$ws.Range("H65").Value = 50000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 50000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 150000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -156864
$ws.Range("H70").Value = 9183.333000000001
$ws.Range("H73").Value = 9183.333000000001
$ws.Range("H134").Value = 199998
$ws.Range("J134").Value = 199998
$ws.Range("L134").Value = 599994
$ws.Range("N134").Value = -605064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 100
$ws.Range("K3").Value = 100
$ws.Range("M3").Value = 12
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 100
$ws.Range("K15").Value = 100
$ws.Range("M15").Value = 70
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H40").Value = 3683.3635
$ws.Range("I40").Value = 2724.111
$ws.Range("K40").Value = 2724.111
$ws.Range("M40").Value = -2588.111
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H93").Value = 2999
$ws.Range("I93").Value = 2999
$ws.Range("K93").Value = 2999
$ws.Range("M93").Value = -1751
$ws.Range("H136").Value = 2727.7856
$ws.Range("I136").Value = 2599.1667
$ws.Range("J136").Value = 3499.5
$ws.Range("K136").Value = 7797.500100000001
$ws.Range("L136").Value = 10498.5
$ws.Range("M136").Value = -5247.500100000001
$ws.Range("N136").Value = -15598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
